# TC_22.xlsx edit: rename sheet, tweak a comment/number format/text/value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet name: "My Series" -> "Data"
$ws.Name = "Data"

# 2) Cell A11: "Function Description" -> "Function Information"
$ws.Range("A11").Value = "Function Information"

# 3) Cell B21: Kurtosis value tweak
$ws.Range("B21").Value = 0.2499825759175085

# 4) numFmt 166 (applied to column A of the data rows, e.g. A27:A36 style 7
#    uses numFmt 164; B27:B36 uses style 8 -> numFmt 166) "0.000" -> "###0.000"
$ws.Range("B27:B36").NumberFormat = "###0.000"
